$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the material for Byggmakker (row 5): "Tree" -> "Timber"
$ws.Range("B5").Value = "Timber"

# Move selection to A6 (matches the saved selection state in the file)
$ws.Range("A6").Select()
